$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number that needs to move
# forward by one day (45179 -> 45180, i.e. 2023-09-10 -> 2023-09-11)
# for every data row (rows 2 through 261).
$firstRow = 2
$lastRow = 261

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45179) {
        $cell.Value = 45180
    }
}
